$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Tentative Dates" (sheet1)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Fill in "NA" for a few blank Review-comment cells
$ws1.Range("C9").Value = "NA"
$ws1.Range("C10").Value = "NA"
$ws1.Range("C11").Value = "NA"

# Bump Vincy's progress percentage
$ws1.Range("B5").Value = 0.15

# Date First Review was "29/01/2023" (text) -> becomes a real date value 3/3/2023
$ws1.Range("B16").Value = 44988
$ws1.Range("B16").NumberFormat = "mm-dd-yy"

# Date Second Review changes from 22/01/2023 to 25/02/2023
$ws1.Range("B18").Value = "25/02/2023"

# ---------------------------------------------------------------------------
# Sheet "Mock interview Schedules" (sheet2)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("B4").Value = 44962
$ws2.Range("B6").Value = 44962
$ws2.Range("B8").Value = 44975

$ws2.Range("B10").Value = "NA"
$ws2.Range("B11").Value = "NA"

$ws2.Range("B12").Value = 44972
$ws2.Range("B12").NumberFormat = "d-mmm-yy"

# ---------------------------------------------------------------------------
# Sheet "Project-1" (sheet3)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("D6").Value = "Completed"

# ---------------------------------------------------------------------------
# Sheet "Resume_CV Preparation" (sheet5)
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("D6").Value = "Completed"

# ---------------------------------------------------------------------------
# Selections - apply in order so "Tentative Dates" ends up the active tab
# ---------------------------------------------------------------------------
$ws5.Range("C13").Select()
$ws3.Range("A36").Select()
$ws2.Range("B16").Select()
$ws1.Range("D17").Select()
